$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new measurement columns (C, D, E, F) for the existing rows and
# update/add values so the averages in column H recalculate to the
# values described in the commit.

# Row 2: LEDs=0, mA base 5.9, add C2
$ws.Range("C2").Value = 5.61

# Row 3: LEDs=1, mA base 16.37, add C3:F3
$ws.Range("C3").Value = 16.4
$ws.Range("D3").Value = 15.7
$ws.Range("E3").Value = 25.1
$ws.Range("F3").Value = 16

# Row 4: LEDs=2, mA base 26, add C4:F4
$ws.Range("C4").Value = 25.8
$ws.Range("D4").Value = 25.8
$ws.Range("E4").Value = 26
$ws.Range("F4").Value = 25

# Row 5: LEDs=3, mA base 35, add C5:F5
$ws.Range("C5").Value = 34.7
$ws.Range("D5").Value = 34.9
$ws.Range("E5").Value = 35.1
$ws.Range("F5").Value = 34.3

# Row 6: LEDs=4, mA base 44.1, add C6
$ws.Range("C6").Value = 43.6

$wb.RefreshAll()
$excel.CalculateFullRebuild()
